$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dialogs")

# Update existing cell text (rows 6-13) that changed wording
$ws.Range("A6").Value = "book a room in Florida"
$ws.Range("B7").Value = "What date?"
$ws.Range("A8").Value = "02/06/18"
$ws.Range("A10").Value = "10:00 AM"
$ws.Range("A12").Value = "1 hour"

# Append the two new rows (14 and 15) that extend the dialog sequence
$ws.Range("A14").Value = "3 people"
$ws.Range("B15").Value = "Choose a room please."
